$d = $word.ActiveDocument

# The commit swaps the East-Asian font from "DejaVu Sans" to "Tahoma"
# (docDefaults + the Normal/Heading paragraph styles), and also records
# an explicit complex-script font ("DejaVu Sans") on the List, Caption
# and Index paragraph styles, which previously inherited it implicitly.

# Normal style: eastAsia DejaVu Sans -> Tahoma
$d.Styles("Normal").Font.NameFarEast = "Tahoma"

# Heading style: eastAsia DejaVu Sans -> Tahoma
$d.Styles("Heading").Font.NameFarEast = "Tahoma"

# List style: add explicit complex-script (cs) font
$d.Styles("List").Font.NameBi = "DejaVu Sans"

# Caption style: add explicit complex-script (cs) font
$d.Styles("Caption").Font.NameBi = "DejaVu Sans"

# Index style: add explicit complex-script (cs) font
$d.Styles("Index").Font.NameBi = "DejaVu Sans"
